# Apply the "456a3b4" data refresh to 杭州-漫展信息.xlsx
#
# Sheet 1 = 展览 (Exhibitions)
# Sheet 2 = 演出 (Performances)
# Sheet 3 = 本地生活 (Local life)
# Sheet 4 = 全部类型 (All types, aggregates the other three)

$wb = $excel.ActiveWorkbook

$wsExpo   = $wb.Worksheets.Item("展览")
$wsShow   = $wb.Worksheets.Item("演出")
$wsLocal  = $wb.Worksheets.Item("本地生活")
$wsAll    = $wb.Worksheets.Item("全部类型")

# ---------------------------------------------------------------
# Sheet 1: 展览 (展览 tab) - "想去人数" (F) / "最低票价" (G) refresh
# ---------------------------------------------------------------
$wsExpo.Range("F2").Value  = 810
$wsExpo.Range("F3").Value  = 14694
$wsExpo.Range("F4").Value  = 14930
$wsExpo.Range("G4").Value  = "不可售"
$wsExpo.Range("F5").Value  = 6036
$wsExpo.Range("F12").Value = 1953
$wsExpo.Range("F15").Value = 2341
$wsExpo.Range("F18").Value = 3586
$wsExpo.Range("F21").Value = 2636
$wsExpo.Range("F22").Value = 670
$wsExpo.Range("F30").Value = 7425
$wsExpo.Range("F31").Value = 5163
$wsExpo.Range("F32").Value = 327
$wsExpo.Range("F44").Value = 23

# ---------------------------------------------------------------
# Sheet 2: 演出
# ---------------------------------------------------------------
$wsShow.Range("F11").Value = 9

# ---------------------------------------------------------------
# Sheet 3: 本地生活
# ---------------------------------------------------------------
$wsLocal.Range("F2").Value = 7970
$wsLocal.Range("F3").Value = 311
$wsLocal.Range("F4").Value = 1093

# ---------------------------------------------------------------
# Sheet 4: 全部类型 (mirrors the rows above, plus its own row 7
# which switches from the TCD event to the 华盟次元嘉年华 event)
# ---------------------------------------------------------------
$wsAll.Range("F2").Value  = 7970
$wsAll.Range("F3").Value  = 810
$wsAll.Range("F4").Value  = 311
$wsAll.Range("F5").Value  = 1093
$wsAll.Range("F6").Value  = 14694

$wsAll.Range("C7").Value = "杭州·第五届华盟次元嘉年华&周年庆狂欢"
$wsAll.Range("D7").Value = "创意路1号 中国智谷富春园区"
$wsAll.Range("E7").Value = "2024.07.20 10:00-07.21 17:00"
$wsAll.Range("F7").Value = 6036
$wsAll.Range("G7").Value = 68
$wsAll.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=84762"
$wsAll.Range("I7").Value = "//i1.hdslb.com/bfs/openplatform/202407/NSQarDy41720678771123.jpeg"

$wsAll.Range("F17").Value = 3586
$wsAll.Range("F19").Value = 2636
$wsAll.Range("F20").Value = 670
$wsAll.Range("F23").Value = 9
$wsAll.Range("F31").Value = 7425
$wsAll.Range("F32").Value = 5163
$wsAll.Range("F34").Value = 327
$wsAll.Range("F44").Value = 23
